$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 40 and row 41 data (A, B, D, E columns) - internet-computer and first-digital-usd swap ranks
$ws.Range("A40").Value = "internet-computer"
$ws.Range("B40").Value = "internet computer"
$ws.Range("D40").Value = 482076964.6058819
$ws.Range("E40").Value = "https://www.dfinityexplorer.org/#/"

$ws.Range("A41").Value = "first-digital-usd"
$ws.Range("B41").Value = "first digital usd"
$ws.Range("D41").Value = 2588088088.167816
$ws.Range("E41").Value = "Sin datos"

# Update marketCapUsd (column F) values for all rows with refreshed data
$ws.Range("F2").Value = 1641207194690.05
$ws.Range("F3").Value = 220337252974.808
$ws.Range("F4").Value = 144586473851.1635
$ws.Range("F5").Value = 125766898323.9508
$ws.Range("F6").Value = 87094124885.83772
$ws.Range("F7").Value = 64667504451.95541
$ws.Range("F8").Value = 60234813484.98527
$ws.Range("F9").Value = 25346141527.26547
$ws.Range("F10").Value = 23691787244.09586
$ws.Range("F11").Value = 22064228294.13318
$ws.Range("F12").Value = 16968716904.44786
$ws.Range("F13").Value = 10659986060.51545
$ws.Range("F14").Value = 9643303953.745028
$ws.Range("F15").Value = 8967880690.436407
$ws.Range("F16").Value = 8872512198.864807
$ws.Range("F17").Value = 8296295563.330326
$ws.Range("F18").Value = 7992201792.645712
$ws.Range("F19").Value = 7517552122.125976
$ws.Range("F20").Value = 7412444842.72647
$ws.Range("F21").Value = 7217873319.067986
$ws.Range("F22").Value = 6527147267.357023
$ws.Range("F23").Value = 6242045410.163404
$ws.Range("F24").Value = 6140514860.934415
$ws.Range("F25").Value = 6111778234.669128
$ws.Range("F26").Value = 6032029266.771893
$ws.Range("F27").Value = 5481899749.226211
$ws.Range("F28").Value = 5364602397.969604
$ws.Range("F29").Value = 5247501502.827567
$ws.Range("F30").Value = 5237699589.282393
$ws.Range("F31").Value = 4237394546.089977
$ws.Range("F32").Value = 4005469081.424708
$ws.Range("F33").Value = 3752579971.891531
$ws.Range("F34").Value = 3247459454.593124
$ws.Range("F35").Value = 3170749559.195848
$ws.Range("F36").Value = 2973819546.164363
$ws.Range("F37").Value = 2905169116.516918
$ws.Range("F38").Value = 2797454361.216642
$ws.Range("F39").Value = 2667981155.819436
$ws.Range("F40").Value = 2584064946.91233
$ws.Range("F41").Value = 2583280325.353663
$ws.Range("F42").Value = 2536838011.83516
$ws.Range("F43").Value = 2519332462.853556
$ws.Range("F44").Value = 2502723147.548345
$ws.Range("F45").Value = 2132526451.965802
$ws.Range("F46").Value = 2043751583.772997
$ws.Range("F47").Value = 1974168246.487428
$ws.Range("F48").Value = 1960605547.507264
$ws.Range("F49").Value = 1960187670.166676
$ws.Range("F50").Value = 1917352673.083153
$ws.Range("F51").Value = 1911368024.400052
$ws.Range("F52").Value = 1829675548.528838
$ws.Range("F53").Value = 1794346876.414898
$ws.Range("F54").Value = 1687107422.514118
$ws.Range("F55").Value = 1672487758.068358
$ws.Range("F56").Value = 1552633403.390302
$ws.Range("F57").Value = 1547508740.962107
$ws.Range("F58").Value = 1469108246.962049
$ws.Range("F59").Value = 1444885033.765069
$ws.Range("F60").Value = 1367680428.03329
$ws.Range("F61").Value = 1358888696.507099
$ws.Range("F62").Value = 1311417512.585624
$ws.Range("F63").Value = 1226760960.995248
$ws.Range("F64").Value = 1202451897.854906
$ws.Range("F65").Value = 1083740388.234232
$ws.Range("F66").Value = 1065454198.640723
$ws.Range("F67").Value = 1049250389.066715
$ws.Range("F68").Value = 962507580.9806228
$ws.Range("F69").Value = 958571455.1192458
$ws.Range("F70").Value = 953591446.119048
$ws.Range("F71").Value = 924549006.9580034
$ws.Range("F73").Value = 891548861.194496
$ws.Range("F74").Value = 882787928.7271338
$ws.Range("F75").Value = 881726841.4225415
$ws.Range("F76").Value = 876452745.498787
$ws.Range("F77").Value = 875455769.9508352
$ws.Range("F78").Value = 843164617.3574919
$ws.Range("F79").Value = 825090255.409063
$ws.Range("F80").Value = 815748307.4967141
$ws.Range("F81").Value = 801762690.4410471
$ws.Range("F82").Value = 792133992.3010215
$ws.Range("F83").Value = 781963907.7790666
$ws.Range("F84").Value = 767207208.6997688
$ws.Range("F85").Value = 692716043.1744293
$ws.Range("F86").Value = 689328171.0181694
$ws.Range("F87").Value = 682073294.9507409
$ws.Range("F88").Value = 679515144.1914803
$ws.Range("F89").Value = 674628377.5776187
$ws.Range("F90").Value = 664200063.8926384
$ws.Range("F91").Value = 658079319.4493791
$ws.Range("F92").Value = 639076930.0121099
$ws.Range("F93").Value = 635134845.9482838
$ws.Range("F94").Value = 629929521.8492606
$ws.Range("F95").Value = 618044828.8354813
$ws.Range("F96").Value = 604688115.988225
$ws.Range("F97").Value = 604192008.2020433
$ws.Range("F98").Value = 584083082.9055908
$ws.Range("F99").Value = 578906400.6221349
$ws.Range("F100").Value = 569217007.7225569
$ws.Range("F101").Value = 566159850.6601783
